# Apply updated crypto price/volume figures to Sheet1 (columns D and E).
# Each value is written with a leading apostrophe so Excel stores it as literal
# text (matching the source inlineStr cells) instead of auto-coercing numeric-
# looking strings (e.g. "215.78") into actual numbers; the Style reset afterwards
# keeps the cell's style index at the sheet's default (no visual/style change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.773.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.18%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.634.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'215.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.23%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  -0.84%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.19%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  -0.45%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.0634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.37%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'19.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.18%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.35%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +0.27%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.860.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.00%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'1.637.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.42%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "'  +0.57%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  -0.81%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'63.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.16%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'25.794.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.12%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E20").Value = "'  +1.19%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'192.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.69%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.37%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "'  +2.29%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +3.41%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  -0.18%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'142.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.98%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.124"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.12%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  +0.90%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'15.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.17%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  +0.10%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.0494"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.63%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +0.30%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.78%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  -1.12%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -0.31%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +0.28%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'1.131.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.71%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.545"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.45%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -2.23%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -1.02%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.00%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'2.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.42%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  -0.45%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'100.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.05%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.799"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.36%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'1.769.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.08%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'  +0.58%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'55.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.19%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0504"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.23%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  -0.77%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +3.39%  "
$ws.Range("E51").Style = "Normal"
